$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2/B3 held real dates previously (formatted as dates). They now hold
# plain text date-like strings ("2022/01/01", "2023/01/01"), so reset any
# existing number formatting/font quirks on those two cells first, then
# apply a Text format before writing the values so Excel stores them
# verbatim as strings rather than re-parsing them back into serials.
$ws.Range("B2:B3").ClearFormats()
$ws.Range("B2:B3").NumberFormat = "@"

$ws.Range("B2").Value = "2022/01/01"
$ws.Range("B3").Value = "2023/01/01"

# Update the numeric parameter cells (bounding box changed from Kenya to Namibia)
$ws.Range("B4").Value = 11
$ws.Range("B5").Value = 26
$ws.Range("B6").Value = -29
$ws.Range("B7").Value = -16

# Update the filename text
$ws.Range("B8").Value = "Namibia-2022"

# Update the active selection to C5
$ws.Range("C5").Select()
